$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column
$ws.Range("AF1").Value = """+/- ratio"""

# Add formulas for the new ratio column, rows 2 through 39
for ($r = 2; $r -le 39; $r++) {
    $ws.Range("AF$r").Formula = "=(Tabelle4[[#This Row],[1]]+Tabelle4[[#This Row],[3]]+Tabelle4[[#This Row],[5]]+1)/(Tabelle4[[#This Row],[2]]+Tabelle4[[#This Row],[4]]+Tabelle4[[#This Row],[6]]+1)"
}
